$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mis-recorded Kaplan-Meier rows: wherever time=25 & event=1 was
# accidentally entered, it should actually be time=20 & event=0.
for ($r = 2; $r -le 999; $r++) {
    $t = $ws.Cells.Item($r, 1).Value2
    $e = $ws.Cells.Item($r, 2).Value2
    if ($t -eq 25 -and $e -eq 1) {
        $ws.Cells.Item($r, 1).Value = 20
        $ws.Cells.Item($r, 2).Value = 0
    }
}

# Restore the view to the top of the sheet (author scrolled back up after
# the edit) instead of leaving it parked at the bottom near row 995.
$null = $ws.Range("G98").Select()
